# Apply the "test meta and color" edit:
#  - Remove the two obsolete paper rows (Boerner and Kruger 2009,
#    Cassidy et al. 2017), shifting the remaining papers up.
#  - Fill in the previously-missing Pigment value for the
#    Whiting, MJ. 1999 / Platysaurus broadleyi row.
#  - Add a new data row for Yewers et al. 2016 / Ctenophorus decresii
#    (species, color pattern, "throat color" acts/time column, "score" units).
#  - Leave the selection on F7, matching the post-edit cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that are no longer part of the table.
$ws.Rows("7:8").Delete()

# Fill in the Pigment value that was previously blank for Whiting, MJ. 1999.
$ws.Range("C6").Value = "orange/yellow"

# Populate the new row for Yewers et al. 2016 (species/color/acts-or-time/units).
$ws.Range("B7").Value = "lizard (Ctenophorus decresii)"
$ws.Range("C7").Value = "yellow/orange grey"
$ws.Range("D7").Value = "throat color"
$ws.Range("E7").Value = "score"

# Match the saved selection/active cell from the edit.
$ws.Range("F7").Select()
